# Shift the per-trial data for rows 69-76 nine columns to the right
# (A:Z -> J:AI), mirroring the layout already used by rows 1-68, and
# zero-fill the newly vacated A:I columns.
#
# Columns must be copied from the highest index down to the lowest
# because the source range (A:Z) and destination range (J:AI) overlap
# in columns J:Z - copying low-to-high would clobber data before it's
# read.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 69; $r -le 76; $r++) {
    for ($col = 26; $col -ge 1; $col--) {
        $src = $ws.Cells.Item($r, $col)
        $dst = $ws.Cells.Item($r, $col + 9)
        $src.Copy($dst)
    }
    for ($col = 1; $col -le 9; $col++) {
        $cell = $ws.Cells.Item($r, $col)
        $cell.ClearFormats()
        $cell.Value = 0
    }
}

# Update the sheet view: move the active selection to J78 and scroll so
# that O61 becomes the top-left visible cell.
[void]$ws.Range("J78").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 61
$win.ScrollColumn = 15
